$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.225.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.585.24'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('E6').Value = '  -2.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.806.25'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.614.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.213.91'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '211.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.31'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0508'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.286.35'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.602'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.56%  '
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.813'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.14'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.51'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.763'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.719.09'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.73'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('E48').Value = '  -3.69%  '
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0505'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₇0974'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.98%  '
